# Updates the crypto price list: refresh "Price" figures (column D) and
# rotate the Coin/Link/Price/Volume rows 18-24 by one position (a coin
# dropped from the table, shifting everyone up and wrapping "One" to the
# bottom), matching the scraped GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumberFormat is set to "@" (Text) immediately before each price write so
# the numeric-looking string (e.g. "252.07") is stored as text, preserving
# the exact original formatting/trailing zeros instead of being coerced to
# a floating-point number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "252.07"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.60"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.423"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05719"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.414"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.364"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8127"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9396"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1440"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07487"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03157"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03078"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09373"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.728"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001601"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04763"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006422"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005041"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001029"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.711"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.170"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01161"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3305"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1308"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002999"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04031"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006768"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008147"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005758"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4998"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
